$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Plan1")

# Mark item 11 ("Blindar API") as completed ("ok") instead of "Em andamento"
$ws.Range("C12").Value = "ok"

# Leave the final selection on F11, matching the state the sheet was left in
$ws.Range("F11").Select()
